$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: update the letter date "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: split the mailing-address paragraph "999 Story Road, San Jose CA
# 95122" (the standalone paragraph, not the one inside the property-address
# table) into two paragraphs: "999 Story Road" and "San Jose, CA 95122".
# ---------------------------------------------------------------------------
$addrIndex = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "999 Story Road, San Jose CA 95122*") {
        if (-not $para.Range.Information(12)) {
            $addrIndex = $i
            break
        }
    }
}

if ($addrIndex -gt 0) {
    $addrPara = $d.Paragraphs($addrIndex)
    $addrRange = $addrPara.Range
    $addrRange.Text = "999 Story Road"
    $addrRange.InsertParagraphAfter()

    $cityPara = $d.Paragraphs($addrIndex + 1)
    $cityPara.Range.Text = "San Jose, CA 95122"
}

# ---------------------------------------------------------------------------
# Change 3: remove the empty "NoSpacing" paragraph that immediately follows
# the "...Board of Directors" signature line.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*Board of Directors*") {
        $nextPara = $para.Next()
        if ($nextPara -ne $null) {
            $nextText = $nextPara.Range.Text.Trim()
            if ($nextText -eq "") {
                $nextPara.Range.Delete()
            }
        }
        break
    }
}
